# Stable Version 1.1 GM HeatStakes
# Update the merged-group ID suffixes in column A (rows 2-48) on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = "MERGED-GRP1+GRP2-88"
    3  = "MERGED-GRP1+GRP2-80"
    4  = "MERGED-GRP1+GRP2-336"
    5  = "MERGED-GRP1+GRP2-720"
    6  = "MERGED-GRP1+GRP2-168"
    7  = "MERGED-GRP1+GRP2-616"
    8  = "MERGED-GRP1+GRP2-128"
    9  = "MERGED-GRP1+GRP2-448"
    10 = "MERGED-GRP1+GRP2-896"
    11 = "MERGED-GRP1+GRP2-344"
    12 = "MERGED-GRP1+GRP2-664"
    13 = "MERGED-GRP1+GRP2-48"
    14 = "MERGED-GRP1+GRP2-432"
    15 = "MERGED-GRP1+GRP2-816"
    16 = "MERGED-GRP1+GRP2-200"
    17 = "MERGED-GRP1+GRP2-520"
    18 = "MERGED-GRP1+GRP2-32"
    19 = "MERGED-GRP1+GRP2-352"
    20 = "MERGED-GRP1+GRP2-800"
    21 = "MERGED-GRP1+GRP2-120"
    22 = "MERGED-GRP1+GRP2-632"
    23 = "MERGED-GRP1+GRP2-688"
    24 = "MERGED-GRP1+GRP2-136"
    25 = "MERGED-GRP1+GRP2-520"
    26 = "MERGED-GRP1+GRP2-904"
    27 = "MERGED-GRP1+GRP2-352"
    28 = "MERGED-GRP1+GRP2-800"
    29 = "MERGED-GRP1+GRP2-120"
    30 = "MERGED-GRP1+GRP2-568"
    31 = "MERGED-GRP1+GRP2-952"
    32 = "MERGED-GRP1+GRP2-336"
    33 = "MERGED-GRP1+GRP2-720"
    34 = "MERGED-GRP1+GRP2-168"
    35 = "MERGED-GRP1+GRP2-552"
    36 = "MERGED-GRP1+GRP2-872"
    37 = "MERGED-GRP1+GRP2-384"
    38 = "MERGED-GRP1+GRP2-768"
    39 = "MERGED-GRP1+GRP2-152"
    40 = "MERGED-GRP1+GRP2-536"
    41 = "MERGED-GRP1+GRP2-920"
    42 = "MERGED-GRP1+GRP2-368"
    43 = "MERGED-GRP1+GRP2-752"
    44 = "MERGED-GRP1+GRP2-72"
    45 = "MERGED-GRP1+GRP2-584"
    46 = "MERGED-GRP1+GRP2-968"
    47 = "MERGED-GRP1+GRP2-288"
    48 = "MERGED-GRP1+GRP2-736"
}

foreach ($row in $values.Keys) {
    $ws.Range("A$row").Value = $values[$row]
}
